$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("g11.3")

# Row 30: Sergipe, 2025
$ws.Range("C30").Value = 6.709622526766052
$ws.Range("D30").Value = 26.25038458442191
$ws.Range("E30").Value = 67.03999288881207

# Row 59: Nordeste, 2025
$ws.Range("C59").Value = 1.07206591321034
$ws.Range("D59").Value = 38.95581830759156
$ws.Range("E59").Value = 59.97211577919808

# Row 82: Brasil, 2019
$ws.Range("D82").Value = 39.20083455047176

# Row 85: Brasil, 2022
$ws.Range("C85").Value = 2.460947153736296
$ws.Range("E85").Value = 56.4347030603654

# Row 88: Brasil, 2025
$ws.Range("C88").Value = 1.692187681270074
$ws.Range("D88").Value = 37.34177051606286
$ws.Range("E88").Value = 60.96604180266707
